$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 ("Draw rubberband"), which shifts rows 8-44 down to 9-45.
# The newly inserted row inherits the formatting (style) of the row above it (row 7 -> s18).
$ws.Rows.Item(8).Insert()
$ws.Range("B8").Value = "Draw rubberband"

# B9 (was B8 pre-insert): "Mouse-ray-heightmap intersection testing" -> rich text with "Implement" prefix in bold
$c9 = $ws.Range("B9")
$c9.Value = "Implement Mouse-ray-heightmap intersection testing"
$c9.Characters(1, 9).Font.Bold = $true

# B12 (was B11 pre-insert): "Allow creation of objects: {Box, Sphere}" -> "Allow creation of primitive objects: {Box, Sphere}"
$ws.Range("B12").Value = "Allow creation of primitive objects: {Box, Sphere}"

# B14 (was B13 pre-insert): "Allow importing of 3ds models" -> "Allow importing of 3ds models, use Engine's 3DS File loader"
$ws.Range("B14").Value = "Allow importing of 3ds models, use Engine's 3DS File loader"

# B2: "Add SMaterialLayer containg all maps and a blend mask" -> bold the trailing "blend mask"
$c2 = $ws.Range("B2")
$c2.Value = "Add SMaterialLayer containg all maps and a blend mask"
$c2.Characters(45, 10).Font.Bold = $true

# B3: style change - apply the greyish (theme 2, darker tint) font used for de-emphasized text
$c3 = $ws.Range("B3")
$c3.Font.ThemeColor = 4
$c3.Font.TintAndShade = -0.249977111117893

Write-Output "done"
